# "fixed horizontal centering on registers"
#
# The register sheet should be centered horizontally when printed.
# Make sure the existing print scale (97%) is kept explicitly before
# turning on horizontal centering, then flip centering on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.PageSetup.Zoom = 97
$ws.PageSetup.CenterHorizontally = $true

# Update the sheet's remembered active cell/selection.
$ws.Range("F14").Select()
